$wb = $excel.ActiveWorkbook

# --- Sheet: Restricciones_del_lider ---
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
$ws.Range("A2").Value = "-16.45 + x_1 + x_2 + y_1 - 2y_2"
$ws.Range("B2").Value = -23.55
$ws.Range("D2").Value = 0.86

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

$ws.Range("A2").Value = "20.12839819112222 - x_1 + 0.08539807211710101y_1 + 2.0974413899797693y_2"
$ws.Range("B2").Value = -10.12839819112222
$ws.Range("D2").Value = 0.32
$ws.Range("E2").Value = 6.8999999999999995
$ws.Range("F2").Value = 5.8

$ws.Range("A3").Value = "5.159464476972512 - x_2 - 0.991360228489825y_1 + 0.8688325598000715y_2"
$ws.Range("B3").Value = 4.840535523027488
$ws.Range("D3").Value = 0.82
$ws.Range("E3").Value = 9.7
$ws.Range("F3").Value = 0

$ws.Range("A4").Value = "24.469905985957396 + 0.9752945376651194y_1 + 2.25386171605379y_2"
$ws.Range("B4").Value = -34.469905985957396
$ws.Range("D4").Value = 0.26
$ws.Range("E4").Value = 5.6000000000000005
$ws.Range("F4").Value = 8.299999999999999

$ws.Range("A5").Value = "-17.988035658026536 - 0.5423001595735621y_1 - 1.7598040282313723y_2"
$ws.Range("B5").Value = -2.0119643419734636
$ws.Range("D5").Value = 0.85
$ws.Range("E5").Value = 8.4
$ws.Range("F5").Value = 4.1

$ws.Range("A6").Value = "-41.5591050815185 + 2.0179935737236696y_1 + 1.3025824110436748y_2"
$ws.Range("B6").Value = -31.559105081518503
$ws.Range("D6").Value = 0.7
$ws.Range("E6").Value = 9.9
$ws.Range("F6").Value = 8.2

$ws.Range("A7").Value = "-41.785557538974174 - 0.45855051767225985y_1 + 0.47678210163037016y_2"
$ws.Range("B7").Value = -21.785557538974174
$ws.Range("D7").Value = 0.51
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0.2

# --- Sheet: Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = 1.55
$ws.Range("B2").Value = 2.7
$ws.Range("C2").Value = -5.1
$ws.Range("D2").Value = -8.65

# --- Sheet: Vector_bf (lowercase "bf") is worksheet #5 ---
# NOTE: Worksheets.Item(name) lookup is case-insensitive, and this workbook has
# both "Vector_bf" and "Vector_BF" sheets which only differ by case, so we must
# address them by their (1-based) position instead of by name to avoid ambiguity.
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = -26.885768177464936
$ws.Range("A3").Value = -18.928761125568965

# --- Sheet: Vector_BF (uppercase "BF") is worksheet #6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 8.719999999999999
$ws.Range("A3").Value = 12.18
$ws.Range("A4").Value = -14.397516931627774
$ws.Range("A5").Value = -34.254859063011175

# --- Sheet: Vector_Alpha ---
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 2.34
$ws.Range("A3").Value = 2.67
